$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FilesTab (row 4) Cypher query (B4): drop the "File Type" and "Breed" RETURN columns
# per "corrected ICDC Breed 1-14 scripts"
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed  IN ['Mixed Breed']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# Row 4 is word-wrapped (style s="1"); removing two lines shrinks its autofit height
$ws.Rows(4).RowHeight = 217.5

# Selection moves onto the edited cell B4
$ws.Range("B4").Select()
